$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 9, 10, 11: set Started (E) and Completed (F) dates, and Status (G) to "Done"
foreach ($r in 9,10,11) {
    $ws.Cells.Item($r, 5).Value = [DateTime]"2013-02-15"
    $ws.Cells.Item($r, 6).Value = [DateTime]"2013-02-15"
    $ws.Cells.Item($r, 7).Value = "Done"
}

# Row 27: set Completed (F) date, and Status (G) to "Done"
$ws.Cells.Item(27, 6).Value = [DateTime]"2013-02-14"
$ws.Cells.Item(27, 7).Value = "Done"

# Update the sheet view: move selection to D20
$ws.Range("D20").Select()
